$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 159, shifting existing rows 159:257 down to 160:258
$ws.Rows("159:159").Insert()

# Fill the newly inserted row 159 with the new data record
$ws.Range("A159").Value = 6
$ws.Range("B159").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C159").Value = "Metropolitana"
$ws.Range("D159").Value = 44830
$ws.Range("D159").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E159").Value = 13
$ws.Range("F159").Value = 100112022
$ws.Range("G159").Value = "Arveja Verde"
$ws.Range("H159").Value = "Perfection"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 230
$ws.Range("K159").Value = 25000
$ws.Range("L159").Value = 25000
$ws.Range("M159").Value = 25000
$ws.Range("N159").Value = "`$/malla 25 kilos"
$ws.Range("O159").Value = "Provincia de Limarí"
$ws.Range("P159").Value = 1000
$ws.Range("Q159").Value = 25
$ws.Range("R159").Value = "Hortaliza"
